$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new rows at position 13 (pushes the old rows 13.. down to
#    15..), making room for the "docentes responsaveis" data rows.  Every
#    row below 13 keeps its original content/formatting, just shifted down
#    by two -- so from here on we only need to touch the cells whose TEXT
#    actually changed.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Resize(2).Insert()

# ---------------------------------------------------------------------------
# 2. New correct text for the "Objetivos:" row (row 10, unchanged position)
# ---------------------------------------------------------------------------
$objetivos = "Desenvolver os fundamentos da Mecânica Clássica com ênfase no formalismo, suas consequências e aplicações. Ao final do curso, o estudante estará apto a aplicar os diversos formalismos da Mecânica Clássica à descrição do movimento de sistemas de partículas e corpos rígidos com alto grau de complexidade."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# ---------------------------------------------------------------------------
# 3. Fill the two newly inserted rows (13 and 14) with the "Docentes
#    responsaveis" data -- only columns B and C, column A stays empty, just
#    like every other "no A-label" row in this sheet (Avaliação: style).
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("B13").Font.Bold = $false
$ws.Range("B13").VerticalAlignment = -4160
$ws.Range("B13").WrapText = $true

$ws.Range("C13").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C13").Font.Bold = $false
$ws.Range("C13").VerticalAlignment = -4160
$ws.Range("C13").WrapText = $true
$ws.Range("C13").Font.Color = 255

$ws.Range("B14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("B14").Font.Bold = $false
$ws.Range("B14").VerticalAlignment = -4160
$ws.Range("B14").WrapText = $true

$ws.Range("C14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C14").Font.Bold = $false
$ws.Range("C14").VerticalAlignment = -4160
$ws.Range("C14").WrapText = $true
$ws.Range("C14").Font.Color = 255

# The two new rows got a blank, bold "column A" cell carried over from the
# insert -- drop it completely so the row has no A-cell at all (matches the
# other label-less rows such as "Avaliação:").
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()

# ---------------------------------------------------------------------------
# 4. Row 15 ("Programa resumido:") -- new correct Portuguese summary text
#    (previously row 13, wrongly filled with the activation date)
# ---------------------------------------------------------------------------
$programaResumido = "Revisão da dinâmica newtoniana (“mecânica vetorial”). Estrutura geral da mecânica clássica nas formulações lagrangiana e hamiltoniana. Aplicações à problemas de forças centrais e dinâmica de corpos rígidos. Problemas em referenciais não-inerciais."
$ws.Range("B15").Value = $programaResumido
$ws.Range("C15").Value = $programaResumido
$ws.Rows.Item(15).RowHeight = 60

# Row 16 ("Short syllabus:") keeps its previous (already correct) text, just
# make sure the row height matches the target (60, was 120 before the shift).
$ws.Rows.Item(16).RowHeight = 60

# ---------------------------------------------------------------------------
# 5. Row 17 ("Programa:") -- new correct Portuguese full syllabus text
#    (previously row 15, wrongly filled with the teacher's name)
# ---------------------------------------------------------------------------
$programa = "1. Revisão da Cinemática do ponto: vetores posição, velocidade e aceleração. 2. Forças, força resultante, conservação do momento linear e as Leis de Newton da Dinâmica; 3. Trabalho e energia; energia cinética e energia potencial. Teorema da energia cinética e conservação da energia total. 4. Princípio variacional de Hamilton e a Formulação Lagrangiana da Mecânica: coordenadas generalizadas, Lagrangiana e as equações de Euler-Lagrange. Coordenadas ignoráveis e princípios de conservação. 5. Aplicações: forças centrais; torque e conservação do momento angular; problemas de dois corpos com atração ou repulsão mútua; Gravitação e as leis de Kepler. 6. Referenciais não inerciais linearmente acelerados e em rotação. Forças de inércia: força centrífuga, força de Coriolis. Efeitos das forças de inércia no planeta Terra; Pêndulo de Foucault. 7. Estática e Dinâmica de corpos rígidos; momentos de inércia; torques; Movimento plano de corpos rígidos; teorema dos eixos paralelos . 8. Formulação de Hamilton para a Mecânica Clássica: a Hamiltoniana e as equações de Hamilton"
$ws.Range("B17").Value = $programa
$ws.Range("C17").Value = $programa
$ws.Rows.Item(17).RowHeight = 120

# Row 18 ("Syllabus:") keeps its previous (already correct) text.
$ws.Rows.Item(18).RowHeight = 120

# ---------------------------------------------------------------------------
# 6. Row 20 ("Método:") -- correct evaluation-method text (previously held
#    the teacher's name, which now lives in row 14)
# ---------------------------------------------------------------------------
$metodo = "A avaliação será composta por duas provas escritas (P1 e P2)."
$ws.Range("B20").Value = $metodo
$ws.Range("C20").Value = $metodo

# Row 21 ("Critério:") -- shifted-up content (grading-weight description)
$criterio = "A nota final (NF) será a média ponderada de três provas, P1 (peso 1), P2 (peso 1) e P3 (peso 2)"
$ws.Range("B21").Value = $criterio
$ws.Range("C21").Value = $criterio

# Row 22 ("Norma de recuperação:") -- shifted-up content (recovery exam rule)
$norma = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("B22").Value = $norma
$ws.Range("C22").Value = $norma
$ws.Rows.Item(22).RowHeight = 60

# ---------------------------------------------------------------------------
# 7. Row 23 ("Bibliografia:") -- brand-new bibliography text (formatting was
#    already correct after the shift, only the text itself changed)
# ---------------------------------------------------------------------------
$bibliografia = "TAYLOR, J. R. - Mecânica Clássica, Bookman, 2015. THORNTON, S. T. MARION, J. B. – Dinâmica Clássica de Partículas e Sistemas, tradução da 5ª edição norte-americana, CENGAGE Learning, 2016. F.P. BEER, E.R. JOHNSTON, E. RUSSEL. - Mecânica vetorial para engenheiros: Estática, McGraw Hill. 9a Ed., 2012. BEER, F.P., JOHNSTON Jr., E.R., CLAUSEN, W. E. - Mecânica Vetorial para Engenheiros: Dinâmica, McGraw-Hill. 7ª Ed., 2006. GOLDSTEIN, H.; POOLE, C.; SAFKO, J. – Classical Mechanics, Addison-Wesley Pub. Co. 2013.LEMOS, N. A. – Mecânica Analítica, Livraria da Física. 2007.KOMPANEYETS, A. S. – Theoretical Physics, Peace Publishers. 2012. LANDAU, L. D.; LIFSHITZ, E. M. – Mechanics, Pergamon Press. 1969"
$ws.Range("B23").Value = $bibliografia
$ws.Range("C23").Value = $bibliografia
$ws.Rows.Item(23).RowHeight = 120

# Row 24 ("Requisitos:") keeps its label only -- nothing else to change,
# it already carries no B/C value after the shift.

# ---------------------------------------------------------------------------
# 8. Rows 25-26: requisites text (previously rows 23-24, shifted down by 2,
#    text and formatting unchanged) -- nothing to do here.
# ---------------------------------------------------------------------------
